$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 04:22"

# --- Update Noruega (row 32) ---
$ws.Range("B32").Value = 6605
$ws.Range("C32").Value = 2
$ws.Range("E32").Value = 6439

# --- Update Mexico (row 37) ---
$ws.Range("F37").Value = 207

# --- Swap Martinica / Paraguay order + refresh their stats (rows 123-124) ---
# Row 123 currently holds Martinica; it becomes Paraguay with new numbers.
# Row 124 currently holds Paraguay; it becomes Martinica, keeping Martinica's old numbers.
$ws.Range("A123").Value = "Paraguay"
$ws.Range("B123").Value = 159
$ws.Range("C123").Value = 12
$ws.Range("D123").Value = 22
$ws.Range("E123").Value = 130
$ws.Range("F123").Value = 1
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 7

$ws.Range("A124").Value = "Martinica"
$ws.Range("B124").Value = 157
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 50
$ws.Range("E124").Value = 101
$ws.Range("F124").Value = 19
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 6

# --- Update San Martin (Parte Holandesa) (row 151) ---
$ws.Range("B151").Value = 52
$ws.Range("C151").Value = 2
$ws.Range("E151").Value = 38
